$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update nombre_aides (C) and montant_total (D) columns for rows with new 2020-10-04 data
$updates = @(
    @{Row=2; C=40138; D=57982265},
    @{Row=3; C=95863; D=140455941},
    @{Row=4; C=32674; D=48370057},
    @{Row=5; C=9261; D=13759690},
    @{Row=6; C=2219; D=3295420},
    @{Row=7; C=204; D=301093},
    @{Row=12; C=43449; D=58864943},
    @{Row=13; C=10193; D=14727973},
    @{Row=14; C=27145; D=39780362},
    @{Row=15; C=8647; D=12832865},
    @{Row=20; C=10690; D=14107702},
    @{Row=21; C=14081; D=20309731},
    @{Row=22; C=32991; D=48380529},
    @{Row=23; C=10630; D=15798232},
    @{Row=27; C=12200; D=16251445},
    @{Row=28; C=8140; D=11773261},
    @{Row=29; C=23587; D=34607750},
    @{Row=30; C=8132; D=12090233},
    @{Row=34; C=8730; D=11523228},
    @{Row=35; C=3512; D=5068760},
    @{Row=36; C=8315; D=12147062},
    @{Row=37; C=3332; D=4940961},
    @{Row=41; C=2646; D=3570849},
    @{Row=42; C=18246; D=26346299},
    @{Row=43; C=53535; D=78438127},
    @{Row=44; C=19712; D=29265018},
    @{Row=45; C=5902; D=8782202},
    @{Row=46; C=1342; D=2002644},
    @{Row=50; C=17692; D=23464937},
    @{Row=51; C=2253; D=3269962},
    @{Row=52; C=7591; D=11153270},
    @{Row=53; C=2542; D=3794684},
    @{Row=57; C=7674; D=10558411},
    @{Row=58; C=1397; D=2599823},
    @{Row=59; C=3443; D=6439324},
    @{Row=60; C=1355; D=2536312},
    @{Row=61; C=459; D=862583},
    @{Row=62; C=159; D=313100},
    @{Row=64; C=2145; D=3706042},
    @{Row=65; C=16360; D=23617131},
    @{Row=66; C=46990; D=68703923},
    @{Row=67; C=16405; D=24374837},
    @{Row=68; C=4787; D=7129774},
    @{Row=69; C=1032; D=1534329},
    @{Row=70; C=86; D=126330},
    @{Row=73; C=15748; D=20696746},
    @{Row=74; C=58175; D=84573257},
    @{Row=75; C=160801; D=236722962},
    @{Row=76; C=69002; D=102782310},
    @{Row=77; C=22284; D=33297128},
    @{Row=78; C=5466; D=8161641},
    @{Row=79; C=356; D=529170},
    @{Row=84; C=5; D=7500},
    @{Row=85; C=57236; D=77427543},
    @{Row=86; C=4952; D=7176169},
    @{Row=87; C=12257; D=18000417},
    @{Row=89; C=1413; D=2110611},
    @{Row=90; C=320; D=476512},
    @{Row=93; C=5736; D=7696940},
    @{Row=94; C=1735; D=2501061},
    @{Row=95; C=5582; D=8224617},
    @{Row=96; C=2039; D=3034681},
    @{Row=97; C=742; D=1111960},
    @{Row=101; C=3820; D=5064857},
    @{Row=102; C=810; D=1460491},
    @{Row=103; C=529; D=1006392},
    @{Row=104; C=194; D=363734},
    @{Row=107; C=11454; D=16605826},
    @{Row=108; C=30450; D=44708964},
    @{Row=109; C=10213; D=15183945},
    @{Row=110; C=2820; D=4204315},
    @{Row=114; C=10238; D=13490417},
    @{Row=115; C=32249; D=46479985},
    @{Row=116; C=69187; D=101207680},
    @{Row=117; C=22220; D=33010199},
    @{Row=118; C=6348; D=9451714},
    @{Row=119; C=1215; D=1815092},
    @{Row=124; C=27002; D=35997373},
    @{Row=125; C=38367; D=55324703},
    @{Row=126; C=80927; D=118279020},
    @{Row=127; C=24920; D=36983445},
    @{Row=128; C=6738; D=10013116},
    @{Row=129; C=1381; D=2048740},
    @{Row=133; C=33325; D=44180800},
    @{Row=134; C=14059; D=20343494},
    @{Row=135; C=33810; D=49636809},
    @{Row=136; C=11952; D=17758284},
    @{Row=137; C=3144; D=4685875},
    @{Row=138; C=546; D=812990},
    @{Row=142; C=11313; D=15054214},
    @{Row=143; C=37444; D=54070147},
    @{Row=144; C=86141; D=126158071},
    @{Row=145; C=25639; D=38085217},
    @{Row=146; C=6749; D=10063945},
    @{Row=147; C=1559; D=2317802},
    @{Row=150; C=30748; D=41387651}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
}
